# Refreshes the cryptocurrency ranking table (GitHub Actions scheduled
# update). Columns B (Coin) and C (Link) only change when a coin's rank
# shifted; columns D (Price) and E (Volume 1h) are refreshed for every row.
# A newly-tracked coin ("Frax") enters the list at rank 33 (row 35), which
# pushes every following coin down by one rank; "Elrond" (previously the
# last row) drops off the bottom of the fixed 50-row table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.378.53"
$ws.Range("E2").Value = "  -3.38%  "

# Row 3
$ws.Range("D3").Value = "1.800.13"
$ws.Range("E3").Value = "  -3.03%  "

# Row 4
$ws.Range("D4").Value = "'1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.45%  "

# Row 5
$ws.Range("D5").Value = "'1.008"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.52%  "

# Row 6
$ws.Range("D6").Value = "'308.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.15%  "

# Row 7
$ws.Range("D7").Value = "'0.4512"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.04%  "

# Row 8
$ws.Range("D8").Value = "'0.3641"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.90%  "

# Row 9
$ws.Range("D9").Value = "'0.07084"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.94%  "

# Row 10
$ws.Range("D10").Value = "'0.8666"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.62%  "

# Row 11
$ws.Range("D11").Value = "'0.07772"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.71%  "

# Row 12
$ws.Range("D12").Value = "'19.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.42%  "

# Row 13
$ws.Range("D13").Value = "1.780.64"
$ws.Range("E13").Value = "  -2.39%  "

# Row 14
$ws.Range("D14").Value = "'5.247"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.69%  "

# Row 15
$ws.Range("D15").Value = "'6.305"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.44%  "

# Row 16
$ws.Range("D16").Value = "'85.97"
$ws.Range("D16").Style = "Normal"

# Row 17
$ws.Range("E17").Value = "  +0.45%  "

# Row 18
$ws.Range("D18").Value = "'0.000008531"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.42%  "

# Row 19
$ws.Range("D19").Value = "'1.008"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.55%  "

# Row 20
$ws.Range("D20").Value = "26.426.83"
$ws.Range("E20").Value = "  -3.28%  "

# Row 21
$ws.Range("D21").Value = "'14.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.94%  "

# Row 22
$ws.Range("D22").Value = "'4.953"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.31%  "

# Row 23
$ws.Range("D23").Value = "'10.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.42%  "

# Row 24
$ws.Range("D24").Value = "'1.974"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.79%  "

# Row 25
$ws.Range("D25").Value = "'149.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.67%  "

# Row 26
$ws.Range("D26").Value = "'17.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.37%  "

# Row 27
$ws.Range("D27").Value = "'1.975"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.90%  "

# Row 28
$ws.Range("D28").Value = "'112.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.29%  "

# Row 29
$ws.Range("D29").Value = "'4.844"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.51%  "

# Row 30
$ws.Range("D30").Value = "'0.08625"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.34%  "

# Row 31
$ws.Range("D31").Value = "'3.028"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.88%  "

# Row 32
$ws.Range("D32").Value = "'0.7257"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.09%  "

# Row 33
$ws.Range("D33").Value = "'4.421"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.05%  "

# Row 34
$ws.Range("E34").Value = "  -5.43%  "

# Row 35
$ws.Range("B35").Value = "Frax"
$ws.Range("C35").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D35").Value = "'1.005"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.34%  "

# Row 36
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "'2.541"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.95%  "

# Row 37
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'1.069"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.87%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01911"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.28%  "

# Row 39
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.05050"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.15%  "

# Row 40
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.868"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.11%  "

# Row 41
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'6.972"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.10%  "

# Row 42
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.4892"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.65%  "

# Row 43
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.1563"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.73%  "

# Row 44
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "'8.093"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.95%  "

# Row 45
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'1.009"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.63%  "

# Row 46
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.4597"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.11%  "

# Row 47
$ws.Range("D47").Value = "'9.937"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.16%  "

# Row 48
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'101.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.48%  "

# Row 49
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'1.577"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.07%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.05987"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.86%  "

# Row 51
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'63.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.75%  "
